# Auto-generated: add rows 48-67 (flashcards data) to the active worksheet
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 48
$ws.Range("A48").Value = 47
$ws.Range("B48").Value = 'Conhecimentos Específicos'
$ws.Range("C48").Value = 'Layout e Arranjos Físicos'
$ws.Range("D48").Value = 'Tipos de Arranjo Físico'
$ws.Range("E48").Value = @'
<ul>
	<li>Posicional (posição fixa)</li>
	<li>Funcional (por processo)</li>
	<li>Celular
<ul>
	<li>Por produto</li>
	<li>Por processo</li>
</ul>
</li>
	<li>Por processo (em linha/fluxo)</li>
</ul>
'@
$ws.Range("F48").Value = 0
$ws.Range("G48").Value = 0

# Row 49
$ws.Range("A49").Value = 48
$ws.Range("B49").Value = 'Conhecimentos Específicos'
$ws.Range("C49").Value = 'Layout e Arranjos Físicos'
$ws.Range("D49").Value = @'
<b>Arranjo Físico Posicional</b>
<i>Caracterísitcas</i>
'@
$ws.Range("E49").Value = @'
<ul>
	<li>recurso transformado permanece fixo</li>
	<li>recursos transformadores se movimentam ao seu redor</li>
	<li>produto é muito grande, pesado, frágil ou complexo</li>
	<li>Características detalhadas:
<ul>
	<li>produção por projeto;</li>
	<li>Elevada complexidade de coordenação das atividades;</li>
	<li>Alto custo unitário</li>
	<li>Espaço limitado para alocar recursos;</li>
	<li>Grande dependência de planejamento logístico e sequenciamento de tarefas.</li>
</ul>
</li>
</ul>
'@
$ws.Range("F49").Value = 0
$ws.Range("G49").Value = 0

# Row 50
$ws.Range("A50").Value = 49
$ws.Range("B50").Value = 'Conhecimentos Específicos'
$ws.Range("C50").Value = 'Layout e Arranjos Físicos'
$ws.Range("D50").Value = @'
<b>Arranjo Físico Funcional (Por processo)</b>
<i>Características</i>
'@
$ws.Range("E50").Value = @'
<ul>
	<li>Recursos similares são agrupados em áreas especializadas</li>
	<li>permitindo que produtos ou clientes sigam roteiros <b>distintos</b> de acordo com suas necessidades</li>
	<li>Características:
<ul>
	<li>Alta flexibilidade;</li>
	<li>Processos semelhantes juntos</li>
	<li>Alto nível de estoques em processo;</li>
	<li>Fluxo não linear, com roteiros personalizados.</li>
</ul>
</li>
</ul>
'@
$ws.Range("F50").Value = 0
$ws.Range("G50").Value = 0

# Row 51
$ws.Range("A51").Value = 50
$ws.Range("B51").Value = 'Conhecimentos Específicos'
$ws.Range("C51").Value = 'Layout e Arranjos Físicos'
$ws.Range("D51").Value = @'
<b>Arranjo Físico Celular</b>
<i>Características</i>
'@
$ws.Range("E51").Value = @'
<ul>
	<li>Forma híbrida entre os arranjos <b>funcional </b>e <b>em linha</b></li>
	<li>Recursos necessários são agrupados em células especializadas</li>
	<li>Quando utilizado:
<ul>
	<li>Produção em bateladas/grandes lotes;</li>
	<li>Alta variedade com volume considerável;</li>
	<li>Necessidade de autonomia operacional nas células.</li>
</ul>
</li>
	<li>Comparativo com funcional:
<ul>
	<li>Enquanto no funcional os recursos são separados por tipo, no celular eles são agrupados por produto ou família de produto.</li>
</ul>
</li>
</ul>
'@
$ws.Range("F51").Value = 0
$ws.Range("G51").Value = 0

# Row 52
$ws.Range("A52").Value = 51
$ws.Range("B52").Value = 'Conhecimentos Específicos'
$ws.Range("C52").Value = 'Layout e Arranjos Físicos'
$ws.Range("D52").Value = @'
<b>Arranjo Celular por Produto</b>
<i>Características</i>
'@
$ws.Range("E52").Value = @'
<ul>
	<li>cada célula é organizada para atender a um produto específico ou uma família de produtos semelhantes</li>
	<li>O foco está no resultado final da célula, ou seja, no produto a ser entregue</li>
	<li>Os recursos dentro da célula são agrupados com o objetivo de produzir aquele produto do início ao fim (ou grande parte dele).</li>
<li>Características:
<ul>
	<li>Ideal para variedade moderada com repetição;</li>
	<li>Processos tendem a ser sequenciais;</li>
	<li>Pode se assemelhar ao arranjo por produto, mas com flexibilidade.</li>
</ul>
</li>
</ul>
'@
$ws.Range("F52").Value = 0
$ws.Range("G52").Value = 0

# Row 53
$ws.Range("A53").Value = 52
$ws.Range("B53").Value = 'Conhecimentos Específicos'
$ws.Range("C53").Value = 'Layout e Arranjos Físicos'
$ws.Range("D53").Value = @'
<b>Arranjo Celular por Processo</b>
<i>Características</i>
'@
$ws.Range("E53").Value = @'
<ul>
	<li>a célula é formada por processos similares ou complementares</li>
	<li>foco está na especialização funcional dentro da célula</li>
	<li>produtos passando por diferentes células de acordo com o roteiro necessário.</li>
	<li>Características:
<ul>
	<li>Produtos diferentes podem passar pelas mesmas células;</li>
	<li>Mais próximo do arranjo funcional, porém com integração interna maior.</li>
</ul>
</li>
</ul>

'@
$ws.Range("F53").Value = 0
$ws.Range("G53").Value = 0

# Row 54
$ws.Range("A54").Value = 53
$ws.Range("B54").Value = 'Conhecimentos Específicos'
$ws.Range("C54").Value = 'Layout e Arranjos Físicos'
$ws.Range("D54").Value = @'
<b>Arranjo Físico por Produto (ou Em linha/Fluxo)</b>
<i>Características</i>
'@
$ws.Range("E54").Value = @'
<ul>
	<li>sequência lógica de transformação do produto</li>
	<li>Ideal para processos repetitivos e com alta demanda.</li>
	<li>Características:
<ul>
	<li>Alto nível de padronização;</li>
	<li>Alto nível de padronização;</li>
	<li>Alta eficiência e produtividade;</li>
	<li>Baixa flexibilidade.</li>
</ul>
</li>
</ul>
'@
$ws.Range("F54").Value = 0
$ws.Range("G54").Value = 0

# Row 55
$ws.Range("A55").Value = 54
$ws.Range("B55").Value = 'Conhecimentos Específicos'
$ws.Range("C55").Value = 'Gestão da Inovação'
$ws.Range("D55").Value = '<b>Processo de Inovação</b>'
$ws.Range("E55").Value = @'
<ol>
	<li>Geração de Ideias</li>
	<li>Avaliação</li>
	<li>Experimentação</li>
	<li>Comercialização</li>
	<li>Acompanhamento</li>
</ol>
'@
$ws.Range("F55").Value = 0
$ws.Range("G55").Value = 0

# Row 56
$ws.Range("A56").Value = 55
$ws.Range("B56").Value = 'Conhecimentos Específicos'
$ws.Range("C56").Value = 'Gestão da Inovação'
$ws.Range("D56").Value = 'Modelos de Inovação:  <b>Inovação aberta</b>'
$ws.Range("E56").Value = @'
<ul>
	<li>cunhado pelo pesquisador Henry Chesbrough, em 2003</li>
	<li>Chesbrough observou que as inovações mais disruptivas aconteciam com a ajuda de colaboradores externos e empresas externas.</li>
	<li>acontece de 3 formas: <ul> <li>Inbound</li> <li>Outbound</li> <li>Coupled</li> </ul></li>
</ul>
'@
$ws.Range("F56").Value = 0
$ws.Range("G56").Value = 0

# Row 57
$ws.Range("A57").Value = 56
$ws.Range("B57").Value = 'Conhecimentos Específicos'
$ws.Range("C57").Value = 'Gestão da Inovação'
$ws.Range("D57").Value = 'Modelos de Inovação: <b>Inovação fechada</b>'
$ws.Range("E57").Value = @'
<ul>
	<li>modelo de inovação que precede a inovação aberta</li>
	<li>inovação que não cruza as paredes do ambiente da empresa</li>
	<li>acontece uma limitação do que a ideia pode se tornar</li>
</ul>
'@
$ws.Range("F57").Value = 0
$ws.Range("G57").Value = 0

# Row 58
$ws.Range("A58").Value = 57
$ws.Range("B58").Value = 'Conhecimentos Específicos'
$ws.Range("C58").Value = 'Gestão da Inovação'
$ws.Range("D58").Value = 'Indicadores de Inovação'
$ws.Range("E58").Value = @'
<ul>
	<li>Redução de Custos</li>
	<li>Pesquisa e Desenvolvimento</li>
	<li>Investimento Médio por Projeto</li>
	<li>ROI (Return on Investment)</li>
	<li>Ideias Geradas</li>
	<li>Ideias por Colaborador</li>
	<li>Tempo de Comercialização</li>
	<li>Projetos em Andamento</li>
	<li>Quantidade de Inovações</li>
	<li>Quantidade de Patentes</li>
	<li>Tempo Gasto por Projeto</li>
	<li>Taxa de Sucesso</li>
</ul>
'@
$ws.Range("F58").Value = 0
$ws.Range("G58").Value = 0

# Row 59
$ws.Range("A59").NumberFormat = "@"
$ws.Range("A59").Value = '58'
$ws.Range("B59").Value = 'Conhecimentos Específicos'
$ws.Range("C59").Value = 'Gestão da Inovação'
$ws.Range("D59").Value = 'Tipos de Inovação'
$ws.Range("E59").Value = @'
<ul>
	<li>Inovação de Processos:</li>
	<li>Inovação de Produto:</li>
	<li>Inovação de Serviços:</li>
	<li>Inovação Organizacional:</li>
	<li>Inovação de Marketing:</li>
	<li>Inovação Incremental:</li>
	<li>Inovação Radical: <ul> <li>Incerteza Técnica:</li> <li>Incertezas de Mercado:</li> <li>Incertezas Organizacionais:</li> <li>Incertezas de Recursos:</li> </ul></li>
	<li>Inovação Disruptiva:</li>
</ul>
'@
$ws.Range("F59").Value = 0
$ws.Range("G59").Value = 0

# Row 60
$ws.Range("A60").Value = 59
$ws.Range("B60").Value = 'Conhecimentos Específicos'
$ws.Range("C60").Value = 'Gestão da Inovação'
$ws.Range("D60").Value = 'Inovação em modelo de negócios'
$ws.Range("E60").Value = @'
baseada no conceito do BMC (Business Model Canvas), criado por Osterwalder e Pigneur (2010)
princípios e dimensões <ul> <li>Criação de Valor: <ul> <li>Atividades-chave:</li> <li>Recursos-chave:</li> <li>Parcerias-chave:</li> </ul></li> <li>Entrega de Valor: <ul> <li>Canais:</li> <li>Segmento de Clientes:</li> <li>Relacionamento com o cliente:</li> </ul></li> <li>Captura de Valor: <ul> <li>Estrutura de Custos:</li> <li>Fontes de Receita:</li> </ul></li> </ul>
inovação que mais afeta a empresa no geral
envolve diversos tipos de outras inovações
quatro áreas para se focar e mudar as coisas <ul> <li>Oferta</li> <li>Consumidores</li> <li>Processos</li> <li>Canais de Entrega</li> </ul>
fatores que levam uma empresa a investir na inovação de Modelo de Negócios <ul> <li>Reação</li> <li>Adaptação</li> <li>Adaptação</li> <li>Proatividade</li> <li></ul>
imprevisível e caótica
atitudes que são levadas em conta na hora de montar e avaliar um novo modelo de negócios: <ul> <li>Atitude de Design</li> <li>Atitude de Decisão.</li> </ul>
Processo de inovação em modelo de negócios (Business Model Generation): <ul> <li>Mobilização</li> <li>Compreensão</li> <li>Design</li> <li>Implementação</li> <li>Gerenciamento</li> </ul>
'@
$ws.Range("F60").Value = 0
$ws.Range("G60").Value = 0

# Row 61
$ws.Range("A61").Value = 60
$ws.Range("B61").Value = 'Conhecimentos Específicos'
$ws.Range("C61").Value = 'Gestão da Inovação'
$ws.Range("D61").Value = 'Cultura Organizacional para a Inovação'
$ws.Range("E61").Value = @'
<b>pilares de uma organização com cultura inovadora:</b>
<ul> <li>Valores</li> <li>Visão</li> <li>Transparência</li> <li>Regras</li> <li>Ambientes</li> <li>Flexibilidade</li> <li>Gestão do Tempo</li> </ul>
'@
$ws.Range("F61").Value = 0
$ws.Range("G61").Value = 0

# Row 62
$ws.Range("A62").Value = 61
$ws.Range("B62").Value = 'Conhecimentos Específicos'
$ws.Range("C62").Value = 'Gestão de Projetos'
$ws.Range("D62").Value = 'Tipos de Projetos'
$ws.Range("E62").Value = @'
<ul>
	<li>social</li>
	<li>pessoal</li>
	<li>cultural</li>
	<li>empresarial</li>
	<li>de pesquisa</li>
</ul>
'@
$ws.Range("F62").Value = 0
$ws.Range("G62").Value = 0

# Row 63
$ws.Range("A63").Value = 62
$ws.Range("B63").Value = 'Conhecimentos Específicos'
$ws.Range("C63").Value = 'Gestão de Projetos'
$ws.Range("D63").Value = 'Etapas/Fases de um projeto'
$ws.Range("E63").Value = @'
<ol>
	<li>Concepção/Iniciação;</li>
	<li>Planejamento;</li>
	<li>Execução;</li>
	<li>Monitoramento e Controle; e</li>
	<li>Encerramento/Conclusão</li>
</ol>
'@
$ws.Range("F63").Value = 0
$ws.Range("G63").Value = 0

# Row 64
$ws.Range("A64").Value = 63
$ws.Range("B64").Value = 'Conhecimentos Específicos'
$ws.Range("C64").Value = 'Gestão de Projetos'
$ws.Range("D64").Value = 'Etapas/Fases de um projeto'
$ws.Range("E64").Value = @'
<ol>
	<li>Concepção/Iniciação;</li>
	<li>Planejamento;</li>
	<li>Execução;</li>
	<li>Monitoramento e Controle; e</li>
	<li>Encerramento/Conclusão</li>
</ol>
'@
$ws.Range("F64").Value = 0
$ws.Range("G64").Value = 0

# Row 65
$ws.Range("A65").Value = 64
$ws.Range("B65").Value = 'Conhecimentos Específicos'
$ws.Range("C65").Value = 'Gestão de Projetos'
$ws.Range("D65").Value = 'Fórmula do IDC'
$ws.Range("E65").Value = 'IDC = Valor agregado/Custo real'
$ws.Range("F65").Value = 0
$ws.Range("G65").Value = 0

# Row 66
$ws.Range("A66").Value = 65
$ws.Range("B66").Value = 'Conhecimentos Específicos'
$ws.Range("C66").Value = 'Gestão de Projetos'
$ws.Range("D66").Value = '<b>PMBOK</b> - <i>áreas de conhecimento</i>'
$ws.Range("E66").Value = @'
Gerenciamento de integração de projetos <ul> <li>Preparação do termo de abertura do Projeto</li> <li>Desenvolvimento do Plano de Gerenciamento</li> <li>Orientação e gerenciamento da execução</li> <li>Monitoramento do trabalho realizado</li> <li>Controle integrado das Mudanças</li> <li>Encerramento do projeto ou de uma fase do projeto</li> </ul>
Gerenciamento de escopo de projetos <ul> <li>Planejamento do gerenciamento de escopo</li> <li>Coleta dos requisitos</li> <li>Definição do escopo de gerenciamento</li> <li>Criação da EAP</li> <li>Validação do Escopo</li> <li>Controle do Escopo</li> </ul>
Gerenciamento de cronograma <ul> <li>Planejamento do Gerenciamento do Cronograma</li> <li>Definição das Atividades</li> <li>Sequenciamento das Atividades</li> <li>Estimativa da duração das Atividades</li> <li>Desenvolvimento do Cronograma</li> <li>Controle do Cronograma</li> </ul>
Gerenciamento de custos <ul> <li>Planejamento de Custos</li> <li>Estimativa dos Custos</li> <li>Controle dos Custos</li> <li>Estudo dos Custos</li> <li>Revisões no orçamento</li> </ul>
Gerenciamento da qualidade <ul> <li>Planejamento da Gestão de Qualidade</li> <li>Garantia da Qualidade</li> <li>Controle de Qualidade</li> </ul>
Gerenciamento da comunicação <ul> <li>Planejamento do Gerenciamento de Comunicação</li> <li>Gerenciamento das Comunicações</li> <li>Controle das Comunicações</li> </ul>
Gerenciamento dos riscos <ul> <li>Planejamento do Gerenciamento de Riscos</li> <li>Identificação dos Riscos</li> <li>Análise Qualitativa</li> <li>Análise Quantitativa</li> <li>Planejamento de resposta</li> <li>Controle de Riscos</li> </ul>
Gerenciamento de aquisições do projeto <ul> <li>Planejamento do Gerenciamento de Aquisições</li> <li>Condução das Aquisições</li> <li>Controle das Aquisições</li> </ul>
Gerenciamento das partes interessadas no projeto <ul> <li>Iniciação</li> <li>Planejamento</li> <li>Execução</li> <li>Controle</li> </ul>

'@
$ws.Range("F66").Value = 0
$ws.Range("G66").Value = 0

# Row 67
$ws.Range("A67").Value = 66
$ws.Range("B67").Value = 'Conhecimentos Específicos'
$ws.Range("C67").Value = 'Gestão de Projetos'
$ws.Range("D67").Value = 'CPM: definições de <b>folga livre</b> e <b>folga total</b>'
$ws.Range("E67").Value = @'
<b>folga livre:</b> É a quantidade de tempo que uma atividade pode atrasar sem atrasar as atividades posteriores;
<b>folga total:</b> E a quantidade de tempo que uma atividade pode atrasar sem atrasar a data final do projeto;

'@
$ws.Range("F67").Value = 0
$ws.Range("G67").Value = 0

